$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A7: change the numeric 6 into a text "6" (force text storage, like typing '6)
$ws.Range("A7").Value = "'6"

# E7: rewrite the HTML list content describing MRP/ERP evolution
$e7Text = @"
<ul>
	<li><b>MRP I:</b>
<ul>
<li>décadas de 1960–1970</li>
	<li>Planejamento de materiais;</li>
	<li>BOM – Bill of Materials, programação de ordens;</li>
	<li>Objetivo central:
<ul>
	<li>minimizar estoques</li>
	<li>garantir disponibilidade de insumos.</li>
</ul>
</li>
</ul>
</li>
	<li><b>MRP II</b>
<ul>
<li>décadas de 1970–1980</li>
	<li>capacidade produtiva, finanças e simulações;</li>
	<li>planejamento financeiro integrado;</li>
	<li>integração interfuncional.</li>
</ul>
</li>
	<li><b>ERP I</b>
<ul>
<li>década de 1990 em diante</li>
	<li>total de todos os departamentos e funções;</li>
	<li>módulos de vendas, distribuição, contabilidade, RH, manutenção e outros;</li>
	<li>ambientes multiempresa e multinacionais.</li>
</ul>
</li>
	<li><b>ERP II</b>
<ul>
	<li>integração entre empresas via internet e cadeia de suprimentos (SCM).</li>
</ul>
</li>
	<li><b>ERPs baseados em nuvem (SaaS)</b>
<ul>
	<li>ampliam escalabilidade e flexibilidade.</li>
</ul>
</li>
</ul>
"@

$ws.Range("E7").Value = $e7Text

# Keep the row height as it was originally (avoid auto row-height side effects
# from the long multi-line text assignment above)
$ws.Rows.Item(7).RowHeight = 15

# G7: erros count changed from 0 to 1
$ws.Range("G7").Value = 1
